$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new weekly rows before the current row 537 ---
# This shifts the existing rows 537:559 down to 540:562, matching the
# "everything moves down by 3" pattern observed in the diff.
$ws.Range("A537:R539").EntireRow.Insert()

# Shared constant values for this market/category block
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100114014
$categoria = "Betarraga"
$variedad  = "Sin especificar"
$unidad    = "`$/unidad"
$origen    = "Región Metropolitana"
$kgUnid    = 1
$clasif    = "Hortaliza"

function Set-Fila {
    param($fila, $fecha, $calidad, $volumen, $pmin, $pmax, $pprom)

    $ws.Cells.Item($fila, 1).Value  = $mercadoId
    $ws.Cells.Item($fila, 2).Value  = $mercado
    $ws.Cells.Item($fila, 3).Value  = $region
    $ws.Cells.Item($fila, 4).Value  = $fecha
    $ws.Cells.Item($fila, 5).Value  = $codreg
    $ws.Cells.Item($fila, 6).Value  = $catId
    $ws.Cells.Item($fila, 7).Value  = $categoria
    $ws.Cells.Item($fila, 8).Value  = $variedad
    $ws.Cells.Item($fila, 9).Value  = $calidad
    $ws.Cells.Item($fila, 10).Value = $volumen
    $ws.Cells.Item($fila, 11).Value = $pmin
    $ws.Cells.Item($fila, 12).Value = $pmax
    $ws.Cells.Item($fila, 13).Value = $pprom
    $ws.Cells.Item($fila, 14).Value = $unidad
    $ws.Cells.Item($fila, 15).Value = $origen
    $ws.Cells.Item($fila, 16).Value = $pprom
    $ws.Cells.Item($fila, 17).Value = $kgUnid
    $ws.Cells.Item($fila, 18).Value = $clasif
}

# New week (fecha 44509) inserted at the top of this sub-block
Set-Fila 537 44509 "Primera" 52000 90 100 94
Set-Fila 538 44509 "Segunda" 44000 80 85  82
Set-Fila 539 44509 "Tercera" 15000 50 50  50

# New week (fecha 44491) appended at the end of this sub-block
Set-Fila 560 44491 "Primera" 46000 90 100 95
Set-Fila 561 44491 "Segunda" 34000 80 85  82
Set-Fila 562 44491 "Tercera" 12000 60 60  60

# Column D carries a datetime display format (s="2") in this sheet; make
# sure the freshly appended rows pick it up too (the inserted rows above
# already inherit it from the Insert() shift).
$ws.Range("D560:D562").NumberFormat = $ws.Range("D559").NumberFormat
